# fix some bugs; add random
# Rewrites the "Processing Time" sheet (sheet1): replaces the 3-job sample
# table with a 10-job table (J0..J9), re-randomised numeric columns, and
# applies centred Calibri formatting to the new data + job-name cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processing Time")

# ---- new data (job name, LH, RH, BOT) -------------------------------------
$rows = @(
  @("J0", 29, 78, 9),
  @("J1", 43, 90, 75),
  @("J2", 91, 85, 39),
  @("J3", 81, 95, 71),
  @(" J4  ", 14, 6, 22),
  @(" J5  ", 84, 2, 52),
  @(" J6  ", 46, 37, 61),
  @("J7  ", 31, 86, 46),
  @("J8  ", 76, 69, 76),
  @("J9  ", 85, 13, 61)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Row 12: bottom spacer cell (style only, no value)
$ws.Cells.Item(12, 1).Value = ""

# ---- formatting -------------------------------------------------------
# Build each distinct format once on a scratch cell, then propagate with
# Copy + PasteSpecial(formats) so every target cell lands on the *same*
# style index instead of each Range.Font.Name call minting a fresh one.

$scratch1 = $ws.Cells.Item(30, 30)
$scratch1.Font.Name = "Calibri"
$scratch1.NumberFormat = "#,##0.00"
$scratch1.HorizontalAlignment = -4108
$scratch1.Copy()
$ws.Range("B2:D11").PasteSpecial(-4122)
$scratch1.Clear()

$scratch2 = $ws.Cells.Item(31, 30)
$scratch2.Font.Name = "Calibri"
$scratch2.HorizontalAlignment = -4108
$scratch2.Copy()
$ws.Range("A6:A12").PasteSpecial(-4122)
$scratch2.Clear()

$excel.CutCopyMode = $false

# ---- view / print state -------------------------------------------------
[void]$ws.Range("N2").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "edit applied"
